$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2001
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H64").Value = 8921.308000000001
$ws.Range("I64").Value = 4999
$ws.Range("J64").Value = 9634.454
$ws.Range("K64").Value = 4999
$ws.Range("L64").Value = 9634.454
$ws.Range("M64").Value = -4751
$ws.Range("N64").Value = -10130.454

$ws.Range("H67").Value = 8921.308000000001
$ws.Range("I67").Value = 4999
$ws.Range("J67").Value = 9634.454
$ws.Range("K67").Value = 4999
$ws.Range("L67").Value = 9634.454
$ws.Range("M67").Value = -4141
$ws.Range("N67").Value = -11350.454

$ws.Range("H69").Value = 9243.529
$ws.Range("I69").Value = 6350
$ws.Range("J69").Value = 9629.333000000001
$ws.Range("K69").Value = 19050
$ws.Range("L69").Value = 28887.999
$ws.Range("M69").Value = -18176
$ws.Range("N69").Value = -30635.999

$ws.Range("H72").Value = 9243.529
$ws.Range("I72").Value = 6350
$ws.Range("J72").Value = 9629.333000000001
$ws.Range("K72").Value = 57150
$ws.Range("L72").Value = 86663.997
$ws.Range("M72").Value = -52782
$ws.Range("N72").Value = -95399.997

$ws.Range("H74").Value = 6809.269
$ws.Range("J74").Value = 9329.429
$ws.Range("L74").Value = 9329.429
$ws.Range("N74").Value = -11201.429

$ws.Range("H77").Value = 6809.269
$ws.Range("J77").Value = 9329.429
$ws.Range("L77").Value = 46647.145
$ws.Range("N77").Value = -56007.145

$ws.Range("H86").Value = 2438.3333
$ws.Range("I86").Value = 2009.2858
$ws.Range("J86").Value = 2813.75
$ws.Range("K86").Value = 2009.2858
$ws.Range("L86").Value = 2813.75
$ws.Range("M86").Value = -886.2858000000001
$ws.Range("N86").Value = -5059.75

$ws.Range("H89").Value = 2438.3333
$ws.Range("I89").Value = 2009.2858
$ws.Range("J89").Value = 2813.75
$ws.Range("K89").Value = 10046.429
$ws.Range("L89").Value = 14068.75
$ws.Range("M89").Value = -4430.429
$ws.Range("N89").Value = -25300.75

$ws.Range("H100").Value = 5191.3335
$ws.Range("I100").Value = 2873.6365
$ws.Range("J100").Value = 8833.429
$ws.Range("K100").Value = 2873.6365
$ws.Range("L100").Value = 8833.429
$ws.Range("M100").Value = -2332.6365
$ws.Range("N100").Value = -9915.429

$ws.Range("H125").Value = 3221.3076
$ws.Range("I125").Value = 1289
$ws.Range("K125").Value = 11601
$ws.Range("M125").Value = -9141

$ws.Range("H137").Value = 3613192.5
$ws.Range("I137").Value = 7678.2
$ws.Range("K137").Value = 23034.6
$ws.Range("M137").Value = -20484.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 208
$ws.Range("I4").Value = 237.75
$ws.Range("K4").Value = 237.75
$ws.Range("M4").Value = -121.75

$ws.Range("H5").Value = 126.22222
$ws.Range("I5").Value = 126.22222
$ws.Range("K5").Value = 126.22222
$ws.Range("M5").Value = -14.22221999999999

$ws.Range("H63").Value = 10242.857
$ws.Range("I63").Value = 1839.3
$ws.Range("K63").Value = 1839.3
$ws.Range("M63").Value = -1153.3

$ws.Range("H66").Value = 10242.857
$ws.Range("I66").Value = 1839.3
$ws.Range("K66").Value = 9196.5
$ws.Range("M66").Value = -5764.5

$ws.Range("H88").Value = 1726.2222
$ws.Range("I88").Value = 1751.4
$ws.Range("J88").Value = 1694.75
$ws.Range("K88").Value = 1751.4
$ws.Range("L88").Value = 1694.75
$ws.Range("M88").Value = -1345.4
$ws.Range("N88").Value = -2506.75

$ws.Range("H91").Value = 1726.2222
$ws.Range("I91").Value = 1751.4
$ws.Range("J91").Value = 1694.75
$ws.Range("K91").Value = 1751.4
$ws.Range("L91").Value = 1694.75
$ws.Range("M91").Value = -347.4000000000001
$ws.Range("N91").Value = -4502.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 126.22222
$ws.Range("I4").Value = 126.22222
$ws.Range("K4").Value = 126.22222
$ws.Range("M4").Value = -11.22221999999999

$ws.Range("H22").Value = 550.1667
$ws.Range("I22").Value = 550.1667
$ws.Range("K22").Value = 550.1667
$ws.Range("M22").Value = -377.1667

$ws.Range("H54").Value = 3733.1667
$ws.Range("I54").Value = 1849.75
$ws.Range("J54").Value = 7500
$ws.Range("K54").Value = 1849.75
$ws.Range("L54").Value = 7500
$ws.Range("M54").Value = -1365.75
$ws.Range("N54").Value = -8468

$ws.Range("H82").Value = 66442.375
$ws.Range("J82").Value = 78047
$ws.Range("L82").Value = 78047
$ws.Range("N82").Value = -78813

$ws.Range("H85").Value = 66442.375
$ws.Range("J85").Value = 78047
$ws.Range("L85").Value = 78047
$ws.Range("N85").Value = -80699

$ws.Range("H86").Value = 4173
$ws.Range("I86").Value = 4006
$ws.Range("J86").Value = 4256.5
$ws.Range("K86").Value = 4006
$ws.Range("L86").Value = 4256.5
$ws.Range("M86").Value = -2883
$ws.Range("N86").Value = -6502.5

$ws.Range("H89").Value = 4173
$ws.Range("I89").Value = 4006
$ws.Range("J89").Value = 4256.5
$ws.Range("K89").Value = 20030
$ws.Range("L89").Value = 21282.5
$ws.Range("M89").Value = -14414
$ws.Range("N89").Value = -32514.5

$ws.Range("H105").Value = 1499.2354
$ws.Range("I105").Value = 1530.5
$ws.Range("K105").Value = 1530.5
$ws.Range("M105").Value = 216.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 177.8077
$ws.Range("I7").Value = 95.166664
$ws.Range("J7").Value = 363.75
$ws.Range("K7").Value = 95.166664
$ws.Range("L7").Value = 363.75
$ws.Range("M7").Value = 17.833336
$ws.Range("N7").Value = -589.75

$ws.Range("H41").Value = 13617.389
$ws.Range("J41").Value = 14937.1875
$ws.Range("L41").Value = 14937.1875
$ws.Range("N41").Value = -15793.1875

$ws.Range("H86").Value = 105107.2
$ws.Range("I86").Value = 4705.8184
$ws.Range("J86").Value = 227820
$ws.Range("K86").Value = 4705.8184
$ws.Range("L86").Value = 227820
$ws.Range("M86").Value = -3582.8184
$ws.Range("N86").Value = -230066

$ws.Range("H89").Value = 105107.2
$ws.Range("I89").Value = 4705.8184
$ws.Range("J89").Value = 227820
$ws.Range("K89").Value = 23529.092
$ws.Range("L89").Value = 1139100
$ws.Range("M89").Value = -17913.092
$ws.Range("N89").Value = -1150332

$ws.Range("H99").Value = 8917.143
$ws.Range("I99").Value = 1857
$ws.Range("J99").Value = 14212.25
$ws.Range("K99").Value = 1857
$ws.Range("L99").Value = 14212.25
$ws.Range("M99").Value = -359
$ws.Range("N99").Value = -17208.25

$ws.Range("H126").Value = 8917.143
$ws.Range("I126").Value = 1857
$ws.Range("J126").Value = 14212.25
$ws.Range("K126").Value = 5571
$ws.Range("L126").Value = 42636.75
$ws.Range("M126").Value = -3101
$ws.Range("N126").Value = -47576.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1532.7059
$ws.Range("I129").Value = 528
$ws.Range("J129").Value = 3374.6667
$ws.Range("K129").Value = 1584
$ws.Range("L129").Value = 10124.0001
$ws.Range("M129").Value = 3416
$ws.Range("N129").Value = -20124.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 640.4761999999999
$ws.Range("I2").Value = 978.2727
$ws.Range("J2").Value = 268.9
$ws.Range("K2").Value = 978.2727
$ws.Range("L2").Value = 268.9
$ws.Range("M2").Value = -865.2727
$ws.Range("N2").Value = -494.9

$ws.Range("H54").Value = 32998.332
$ws.Range("J54").Value = 32998.332
$ws.Range("L54").Value = 32998.332
$ws.Range("N54").Value = -33778.332

$ws.Range("H80").Value = 287934.5
$ws.Range("I80").Value = 430313.75
$ws.Range("J80").Value = 3176
$ws.Range("K80").Value = 430313.75
$ws.Range("L80").Value = 3176
$ws.Range("M80").Value = -429315.75
$ws.Range("N80").Value = -5172

$ws.Range("H83").Value = 287934.5
$ws.Range("I83").Value = 430313.75
$ws.Range("J83").Value = 3176
$ws.Range("K83").Value = 2151568.75
$ws.Range("L83").Value = 15880
$ws.Range("M83").Value = -2146576.75
$ws.Range("N83").Value = -25864

$ws.Range("H132").Value = 24693000
$ws.Range("I132").Value = 34904496
$ws.Range("J132").Value = 15221.917
$ws.Range("K132").Value = 104713488
$ws.Range("L132").Value = 45665.751
$ws.Range("M132").Value = -104710958
$ws.Range("N132").Value = -50725.751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3862.9092
$ws.Range("I68").Value = 4665.3335
$ws.Range("J68").Value = 2900
$ws.Range("K68").Value = 4665.3335
$ws.Range("L68").Value = 2900
$ws.Range("M68").Value = -3916.3335
$ws.Range("N68").Value = -4398

$ws.Range("H71").Value = 3862.9092
$ws.Range("I71").Value = 4665.3335
$ws.Range("J71").Value = 2900
$ws.Range("K71").Value = 23326.6675
$ws.Range("L71").Value = 14500
$ws.Range("M71").Value = -19582.6675
$ws.Range("N71").Value = -21988

$ws.Range("H122").Value = 4764.696
$ws.Range("I122").Value = 4478.316
$ws.Range("J122").Value = 6125
$ws.Range("K122").Value = 13434.948
$ws.Range("L122").Value = 18375
$ws.Range("M122").Value = -10984.948
$ws.Range("N122").Value = -23275

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 302366.66
$ws.Range("I62").Value = 451550
$ws.Range("K62").Value = 451550
$ws.Range("M62").Value = -450926

$ws.Range("H65").Value = 302366.66
$ws.Range("I65").Value = 451550
$ws.Range("K65").Value = 2257750
$ws.Range("M65").Value = -2254630

$ws.Range("H81").Value = 10333.333
$ws.Range("I81").Value = 7500
$ws.Range("K81").Value = 15000
$ws.Range("M81").Value = -13939

$ws.Range("H84").Value = 10333.333
$ws.Range("I84").Value = 7500
$ws.Range("K84").Value = 75000
$ws.Range("M84").Value = -69696

$ws.Range("H123").Value = 74995.5
$ws.Range("J123").Value = 74995.5
$ws.Range("L123").Value = 74995.5
$ws.Range("N123").Value = -84795.5

$ws.Range("H132").Value = 5445468.5
$ws.Range("I132").Value = 7743284.5
$ws.Range("J132").Value = 14267.272
$ws.Range("K132").Value = 23229853.5
$ws.Range("L132").Value = 42801.81600000001
$ws.Range("M132").Value = -23227323.5
$ws.Range("N132").Value = -47861.81600000001
